# Auto-generated Excel COM-interop script to apply the Mateus_Profits.xlsx diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific
# rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets,
# matching refreshed market-board data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 129.875
$ws.Range("I5").Value = 85.666664
$ws.Range("J5").Value = 262.5
$ws.Range("K5").Value = 85.666664
$ws.Range("L5").Value = 262.5
$ws.Range("M5").Value = 29.333336
$ws.Range("N5").Value = -492.5
$ws.Range("H8").Value = 92.8
$ws.Range("I8").Value = 92.8
$ws.Range("K8").Value = 278.4
$ws.Range("M8").Value = -139.4
$ws.Range("H17").Value = 4349427
$ws.Range("J17").Value = 4349427
$ws.Range("L17").Value = 13048281
$ws.Range("N17").Value = -13048617
$ws.Range("H19").Value = 1839
$ws.Range("I19").Value = 2073.75
$ws.Range("J19").Value = 1786.8334
$ws.Range("K19").Value = 2073.75
$ws.Range("L19").Value = 1786.8334
$ws.Range("M19").Value = -1898.75
$ws.Range("N19").Value = -2136.8334
$ws.Range("H38").Value = 372.33334
$ws.Range("I38").Value = 372.33334
$ws.Range("K38").Value = 1117.00002
$ws.Range("M38").Value = -745.0000199999999
$ws.Range("H39").Value = 925.5454999999999
$ws.Range("I39").Value = 230.2
$ws.Range("J39").Value = 1505
$ws.Range("K39").Value = 690.5999999999999
$ws.Range("L39").Value = 4515
$ws.Range("M39").Value = -394.5999999999999
$ws.Range("N39").Value = -5107
$ws.Range("H40").Value = 5868.222
$ws.Range("I40").Value = 3962.8
$ws.Range("J40").Value = 8250
$ws.Range("K40").Value = 3962.8
$ws.Range("L40").Value = 8250
$ws.Range("M40").Value = -3787.8
$ws.Range("N40").Value = -8600
$ws.Range("H92").Value = 1172.6364
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -12496
$ws.Range("H112").Value = 3534.2856
$ws.Range("J112").Value = 3498.4614
$ws.Range("L112").Value = 10495.3842
$ws.Range("N112").Value = -12711.3842
$ws.Range("H113").Value = 4745
$ws.Range("I113").Value = 4745
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4745
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1491
$ws.Range("N113").Value = ""
$ws.Range("H116").Value = 4629.769
$ws.Range("I116").Value = 4021.3333
$ws.Range("K116").Value = 4021.3333
$ws.Range("M116").Value = -579.3332999999998
$ws.Range("H132").Value = 3347.8865
$ws.Range("I132").Value = 789.6579
$ws.Range("K132").Value = 2368.9737
$ws.Range("M132").Value = 161.0263
$ws.Range("H133").Value = 78249.5
$ws.Range("J133").Value = 78249.5
$ws.Range("L133").Value = 78249.5
$ws.Range("N133").Value = -88369.5
$ws.Range("H136").Value = 153936.67
$ws.Range("J136").Value = 153936.67
$ws.Range("L136").Value = 153936.67
$ws.Range("N136").Value = -164136.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1666.3334
$ws.Range("I30").Value = 499.5
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 499.5
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = -349.5
$ws.Range("N30").Value = -4300
$ws.Range("H32").Value = 1241.5696
$ws.Range("I32").Value = 1271.5264
$ws.Range("K32").Value = 1271.5264
$ws.Range("M32").Value = -984.5264
$ws.Range("H61").Value = 15888.417
$ws.Range("I61").Value = 19416
$ws.Range("K61").Value = 19416
$ws.Range("M61").Value = -19204
$ws.Range("H132").Value = 4910
$ws.Range("I132").Value = 4406.7856
$ws.Range("K132").Value = 13220.3568
$ws.Range("M132").Value = -10690.3568
$ws.Range("H136").Value = 15888.417
$ws.Range("I136").Value = 19416
$ws.Range("K136").Value = 58248
$ws.Range("M136").Value = -55698

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4757.7915
$ws.Range("I99").Value = 3480.4375
$ws.Range("J99").Value = 7312.5
$ws.Range("K99").Value = 3480.4375
$ws.Range("L99").Value = 7312.5
$ws.Range("M99").Value = -1982.4375
$ws.Range("N99").Value = -10308.5
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""
$ws.Range("H134").Value = 4619.8438
$ws.Range("I134").Value = 4691.4517
$ws.Range("K134").Value = 14074.3551
$ws.Range("M134").Value = -11539.3551

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 194.61539
$ws.Range("I7").Value = 300
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -187
$ws.Range("H31").Value = 3256.2546
$ws.Range("I31").Value = 1905.909
$ws.Range("J31").Value = 8657.637000000001
$ws.Range("K31").Value = 1905.909
$ws.Range("L31").Value = 8657.637000000001
$ws.Range("M31").Value = -1610.909
$ws.Range("N31").Value = -9247.637000000001
$ws.Range("H34").Value = 3256.2546
$ws.Range("I34").Value = 1905.909
$ws.Range("J34").Value = 8657.637000000001
$ws.Range("K34").Value = 1905.909
$ws.Range("L34").Value = 8657.637000000001
$ws.Range("M34").Value = -1703.909
$ws.Range("N34").Value = -9061.637000000001
$ws.Range("H58").Value = 4254.816
$ws.Range("I58").Value = 2742.1482
$ws.Range("K58").Value = 2742.1482
$ws.Range("M58").Value = -2539.1482
$ws.Range("H62").Value = 7138.625
$ws.Range("I62").Value = 7100.6
$ws.Range("K62").Value = 7100.6
$ws.Range("M62").Value = -6476.6
$ws.Range("H65").Value = 7138.625
$ws.Range("I65").Value = 7100.6
$ws.Range("K65").Value = 35503
$ws.Range("M65").Value = -32383
$ws.Range("H99").Value = 3644.9
$ws.Range("I99").Value = 3449.875
$ws.Range("K99").Value = 3449.875
$ws.Range("M99").Value = -1951.875
$ws.Range("H126").Value = 3644.9
$ws.Range("I126").Value = 3449.875
$ws.Range("K126").Value = 10349.625
$ws.Range("M126").Value = -7879.625
$ws.Range("H136").Value = 4254.816
$ws.Range("I136").Value = 2742.1482
$ws.Range("K136").Value = 8226.444600000001
$ws.Range("M136").Value = -5676.444600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1330.5
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 1477.7142
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 4433.142599999999
$ws.Range("M39").Value = -606
$ws.Range("N39").Value = -5021.142599999999
$ws.Range("H103").Value = 299
$ws.Range("I103").Value = 299
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 897
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -18
$ws.Range("N103").Value = ""
$ws.Range("H106").Value = 11000
$ws.Range("J106").Value = 11000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -34892
$ws.Range("H117").Value = 1250845.1
$ws.Range("J117").Value = 1429408.8
$ws.Range("L117").Value = 4288226.4
$ws.Range("N117").Value = -4295110.4
$ws.Range("H129").Value = 20259386
$ws.Range("I129").Value = 47763516
$ws.Range("J129").Value = 1006496
$ws.Range("K129").Value = 143290548
$ws.Range("L129").Value = 3019488
$ws.Range("M129").Value = -143285548
$ws.Range("N129").Value = -3029488
$ws.Range("H136").Value = 630
$ws.Range("I136").Value = 630
$ws.Range("K136").Value = 1890
$ws.Range("M136").Value = 3210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3949.9092
$ws.Range("I102").Value = 3216.5557
$ws.Range("K102").Value = 3216.5557
$ws.Range("M102").Value = -1594.5557
$ws.Range("H132").Value = 6324.8
$ws.Range("I132").Value = 3875
$ws.Range("K132").Value = 11625
$ws.Range("M132").Value = -9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2247
$ws.Range("I40").Value = 2247
$ws.Range("K40").Value = 2247
$ws.Range("M40").Value = -2111
$ws.Range("H46").Value = 11465.968
$ws.Range("I46").Value = 2824.3635
$ws.Range("J46").Value = 16218.85
$ws.Range("K46").Value = 2824.3635
$ws.Range("L46").Value = 16218.85
$ws.Range("M46").Value = -2636.3635
$ws.Range("N46").Value = -16594.85
$ws.Range("H61").Value = 252151
$ws.Range("I61").Value = 252151
$ws.Range("K61").Value = 252151
$ws.Range("M61").Value = -251949
$ws.Range("H113").Value = 252151
$ws.Range("I113").Value = 252151
$ws.Range("K113").Value = 252151
$ws.Range("M113").Value = -249981
$ws.Range("H122").Value = 3284
$ws.Range("I122").Value = 3448.3333
$ws.Range("K122").Value = 10344.9999
$ws.Range("M122").Value = -7894.999899999999
$ws.Range("H136").Value = 6233.0625
$ws.Range("I136").Value = 5552.5713
$ws.Range("K136").Value = 16657.7139
$ws.Range("M136").Value = -14107.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9495
$ws.Range("J15").Value = 9495
$ws.Range("L15").Value = 9495
$ws.Range("N15").Value = -10071
$ws.Range("H100").Value = 1421.9375
$ws.Range("I100").Value = 1391.4166
$ws.Range("K100").Value = 2782.8332
$ws.Range("M100").Value = -2241.8332
$ws.Range("H113").Value = 910.2
$ws.Range("I113").Value = 981.5
$ws.Range("J113").Value = 625
$ws.Range("K113").Value = 2944.5
$ws.Range("L113").Value = 1875
$ws.Range("M113").Value = -774.5
$ws.Range("N113").Value = -6215
$ws.Range("H122").Value = 2995
$ws.Range("I122").Value = 1743.3334
$ws.Range("K122").Value = 5230.0002
$ws.Range("M122").Value = -2780.0002
$ws.Range("H126").Value = 3976.0435
$ws.Range("I126").Value = 3592.2
$ws.Range("K126").Value = 10776.6
$ws.Range("M126").Value = -8306.599999999999

Write-Host "Applied Mateus_Profits.xlsx market data refresh."
